$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.754.27"
$ws.Range("E2").Value = "  +3.93%  "
$ws.Range("D3").Value = "1.870.27"
$ws.Range("E3").Value = "  +2.90%  "
$ws.Range("D5").Value = "'277.33"
$ws.Range("E5").Value = "  -0.34%  "
$ws.Range("D7").Value = "'0.5280"
$ws.Range("E7").Value = "  +3.65%  "
$ws.Range("D8").Value = "'0.3417"
$ws.Range("E8").Value = "  -3.33%  "
$ws.Range("D9").Value = "'0.06934"
$ws.Range("E9").Value = "  +4.04%  "
$ws.Range("D10").Value = "'20.02"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("D11").Value = "'0.8027"
$ws.Range("E11").Value = "  -3.00%  "
$ws.Range("D12").Value = "'0.07751"
$ws.Range("E12").Value = "  -2.01%  "
$ws.Range("D13").Value = "1.882.39"
$ws.Range("E13").Value = "  +3.64%  "
$ws.Range("D14").Value = "'90.27"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").Value = "'5.173"
$ws.Range("E15").Value = "  +1.90%  "
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("E19").Value = "  +0.07%  "
$ws.Range("D20").Value = "26.794.27"
$ws.Range("E20").Value = "  +3.92%  "
$ws.Range("D21").Value = "2.116.41"
$ws.Range("E21").Value = "  +3.93%  "
$ws.Range("D22").Value = "'4.749"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").Value = "'6.175"
$ws.Range("E24").Value = "  +1.07%  "
$ws.Range("D25").Value = "'2.368"
$ws.Range("E25").Value = "  +6.34%  "
$ws.Range("D26").Value = "'146.31"
$ws.Range("E26").Value = "  +3.05%  "
$ws.Range("D27").Value = "'17.32"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("D29").Value = "'113.21"
$ws.Range("E29").Value = "  +3.52%  "
$ws.Range("D30").Value = "'4.336"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").Value = "'4.339"
$ws.Range("E31").Value = "  +2.34%  "
$ws.Range("D32").Value = "'0.08902"
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("D33").Value = "'0.04929"
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("D34").Value = "'1.164"
$ws.Range("E34").Value = "  +2.12%  "
$ws.Range("D35").Value = "'0.7291"
$ws.Range("D36").Value = "'2.880"
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "'3.266"
$ws.Range("E37").Value = "  +4.01%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").Value = "'2.320"
$ws.Range("E39").Value = "  -2.56%  "
$ws.Range("D40").Value = "'0.5140"
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("D41").Value = "'0.9467"
$ws.Range("E41").Value = "  -1.97%  "
$ws.Range("D42").Value = "'116.12"
$ws.Range("E42").Value = "  +4.35%  "
$ws.Range("D43").Value = "'6.157"
$ws.Range("E43").Value = "  -1.02%  "
$ws.Range("D44").Value = "'8.101"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "'0.4461"
$ws.Range("E46").Value = "  -2.15%  "
$ws.Range("E47").Value = "  -2.08%  "
$ws.Range("D48").Value = "'0.06083"
$ws.Range("E48").Value = "  +4.27%  "
$ws.Range("D49").Value = "'36.33"
$ws.Range("E49").Value = "  -0.53%  "
$ws.Range("D50").Value = "'9.262"
$ws.Range("E50").Value = "  +0.68%  "
$ws.Range("D51").Value = "'1.489"
$ws.Range("E51").Value = "  -0.88%  "
